# Generate Report for Handback
# Updates file identifiers and timestamps across the Overview, zh-cn and de-de
# sheets, keeping cell values and hyperlink display text in sync.

$wb = $excel.ActiveWorkbook

$newName1 = "f2d62a4c-e913-49d1-af2b-2aba803b41bf.md"
$newName2 = "ffffe0181847-0564-4a35-b1c6-304a42343a77.md"

$newZhCnXlf = "f2d62a4c-e913-49d1-af2b-2aba803b41bf.06ad2a0f9e6932d06e5ae424372b6fa2f1326765.zh-cn.xlf"
$newDeDeXlf = "f2d62a4c-e913-49d1-af2b-2aba803b41bf.06ad2a0f9e6932d06e5ae424372b6fa2f1326765.de-de.xlf"

function Set-CellAndHyperlink($ws, $addr, $value) {
    $ws.Range($addr).Value = $value
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq ('$' + $addr.Substring(0,1) + '$' + $addr.Substring(1))) {
            $hl.TextToDisplay = $value
        }
    }
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newName1
Set-CellAndHyperlink $wsOverview "B2" ("e2e\" + $newName1)
$wsOverview.Range("G2").Value = "2016-08-17 13:00:24"

$wsOverview.Range("A3").Value = $newName2
Set-CellAndHyperlink $wsOverview "B3" ("e2e\" + $newName2)
$wsOverview.Range("G3").Value = "2016-08-17 13:00:24"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-CellAndHyperlink $wsZhCn "A2" $newName1
$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("H2").Value = "2016-08-17 13:00:03"
Set-CellAndHyperlink $wsZhCn "I2" $newName1
$wsZhCn.Range("J2").Value = $newZhCnXlf
$wsZhCn.Range("K2").Value = "2016-08-17 13:00:48"

Set-CellAndHyperlink $wsZhCn "A3" $newName2
$wsZhCn.Range("G3").Value = $newZhCnXlf
$wsZhCn.Range("H3").Value = "2016-08-17 13:00:03"
Set-CellAndHyperlink $wsZhCn "I3" $newName2
$wsZhCn.Range("J3").Value = $newZhCnXlf
$wsZhCn.Range("K3").Value = "2016-08-17 13:00:48"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-CellAndHyperlink $wsDeDe "A2" $newName1
$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("H2").Value = "2016-08-17 13:00:24"
Set-CellAndHyperlink $wsDeDe "I2" $newName1
$wsDeDe.Range("J2").Value = $newDeDeXlf
$wsDeDe.Range("K2").Value = "2016-08-17 13:00:55"

Set-CellAndHyperlink $wsDeDe "A3" $newName2
$wsDeDe.Range("G3").Value = $newDeDeXlf
$wsDeDe.Range("H3").Value = "2016-08-17 13:00:24"
Set-CellAndHyperlink $wsDeDe "I3" $newName2
$wsDeDe.Range("J3").Value = $newDeDeXlf
$wsDeDe.Range("K3").Value = "2016-08-17 13:00:55"
